# Helper: write a value into a cell as literal TEXT (keeps numeric-looking
# strings such as "0.33" or codes with leading zeros such as "014462" from
# being auto-coerced into numbers), while leaving the cell's style back at
# the default ("Normal") once the value has been committed.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q2" sheet, placing the copy right after it,
#    so the original Q2 data is safely preserved on the duplicate sheet.
#    The (now redundant) original tab will be re-purposed below to hold the
#    brand-new "2022-Q3" data, which lets it keep its original sheetId/rId.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy([System.Reflection.Missing]::Value, $q2)

$q2.Name = "2022-Q3"

$q2dup = $wb.Worksheets.Item("2022-Q2 (2)")
$q2dup.Name = "2022-Q2"

# ---------------------------------------------------------------------------
# 2. Replace the contents of the (renamed) "2022-Q3" sheet with the new
#    quarterly fund-holding data.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Cells.Clear()

$summary = $wb.Worksheets.Item("总计")

# Header row - copy the header cell formatting used on the "总计" sheet.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Row 2: 光大保德信汇佳混合A
$q3.Range("A2").Value = 0
$summary.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)
Set-TextCell $q3.Range("B2") "014462"
Set-TextCell $q3.Range("C2") "光大保德信汇佳混合A"
Set-TextCell $q3.Range("D2") "0.33"
Set-TextCell $q3.Range("E2") "43.38"
Set-TextCell $q3.Range("F2") "2.64"
Set-TextCell $q3.Range("G2") "0.0087"
$q3.Range("H2").Value = 7

# Row 3: 光大保德信汇佳混合C
$q3.Range("A3").Value = 1
$summary.Range("A2").Copy()
$q3.Range("A3").PasteSpecial(-4122)
Set-TextCell $q3.Range("B3") "014463"
Set-TextCell $q3.Range("C3") "光大保德信汇佳混合C"
Set-TextCell $q3.Range("D3") "0.03"
Set-TextCell $q3.Range("E3") "43.38"
Set-TextCell $q3.Range("F3") "2.64"
Set-TextCell $q3.Range("G3") "0.0008"
$q3.Range("H3").Value = 7

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: the existing "2022-Q2" row becomes
#    the "2022-Q3" row (with the new totals), and a fresh row is appended
#    below it holding the original "2022-Q2" totals.
# ---------------------------------------------------------------------------
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.01

$summary.Range("A3").Value = 1
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.44
